$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 448.42856
$ws.Range("I28").Value = 448.42856
$ws.Range("K28").Value = 448.42856
$ws.Range("M28").Value = 36.57144

$ws.Range("H38").Value = 446.42856
$ws.Range("J38").Value = 949.5
$ws.Range("L38").Value = 2848.5
$ws.Range("N38").Value = -3592.5

$ws.Range("H39").Value = 288.27274
$ws.Range("I39").Value = 143.25
$ws.Range("J39").Value = 675
$ws.Range("K39").Value = 429.75
$ws.Range("L39").Value = 2025
$ws.Range("M39").Value = -133.75
$ws.Range("N39").Value = -2617

$ws.Range("H62").Value = 7298.6665
$ws.Range("I62").Value = 7298.6665
$ws.Range("K62").Value = 7298.6665
$ws.Range("M62").Value = -6674.6665

$ws.Range("H65").Value = 7298.6665
$ws.Range("I65").Value = 7298.6665
$ws.Range("K65").Value = 36493.3325
$ws.Range("M65").Value = -33373.3325

$ws.Range("H98").Value = 1628.75
$ws.Range("I98").Value = 1628.75
$ws.Range("K98").Value = 1628.75
$ws.Range("M98").Value = -130.75

$ws.Range("H122").Value = 1628.75
$ws.Range("I122").Value = 1628.75
$ws.Range("K122").Value = 4886.25
$ws.Range("M122").Value = -2436.25

$ws.Range("H132").Value = 6084.8335
$ws.Range("I132").Value = 7081.8
$ws.Range("J132").Value = 1100
$ws.Range("K132").Value = 21245.4
$ws.Range("L132").Value = 3300
$ws.Range("M132").Value = -18715.4
$ws.Range("N132").Value = -8360

$ws.Range("H135").Value = 2767
$ws.Range("I135").Value = 2723.7
$ws.Range("K135").Value = 24513.3
$ws.Range("M135").Value = -21978.3

$ws.Range("H138").Value = 4625.222
$ws.Range("J138").Value = 4625.222
$ws.Range("L138").Value = 13875.666
$ws.Range("N138").Value = -24155.666

$ws.Range("H141").Value = 20418.6
$ws.Range("I141").Value = 20418.6
$ws.Range("K141").Value = 61255.8
$ws.Range("M141").Value = -56075.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H113").Value = 23333
$ws.Range("J113").Value = 23333
$ws.Range("L113").Value = 23333
$ws.Range("N113").Value = -32011

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1369
$ws.Range("I11").Value = 186.33333
$ws.Range("J11").Value = 2551.6667
$ws.Range("K11").Value = 186.33333
$ws.Range("L11").Value = 2551.6667
$ws.Range("M11").Value = -46.33332999999999
$ws.Range("N11").Value = -2831.6667

$ws.Range("H26").Value = 39999
$ws.Range("I26").Value = 39999
$ws.Range("K26").Value = 39999
$ws.Range("M26").Value = -39707

$ws.Range("H86").Value = 1040.3334
$ws.Range("I86").Value = 749.5
$ws.Range("K86").Value = 749.5
$ws.Range("M86").Value = 373.5

$ws.Range("H89").Value = 1040.3334
$ws.Range("I89").Value = 749.5
$ws.Range("K89").Value = 3747.5
$ws.Range("M89").Value = 1868.5

$ws.Range("H105").Value = 1498.75
$ws.Range("I105").Value = 1150
$ws.Range("K105").Value = 1150
$ws.Range("M105").Value = 597

$ws.Range("H107").Value = 1326.2858
$ws.Range("I107").Value = 1346.25
$ws.Range("K107").Value = 1346.25
$ws.Range("M107").Value = 573.75

$ws.Range("H134").Value = 2800
$ws.Range("I134").Value = 1750
$ws.Range("K134").Value = 5250
$ws.Range("M134").Value = -2715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 7750
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 15000
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = -361
$ws.Range("N10").Value = -15278

$ws.Range("H16").Value = 3665
$ws.Range("I16").Value = 1584
$ws.Range("K16").Value = 1584
$ws.Range("M16").Value = -1297

$ws.Range("H17").Value = 19999.334
$ws.Range("I17").Value = 19999.334
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 19999.334
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -19825.334
$ws.Range("N17").ClearContents()

$ws.Range("H22").Value = 599.6667
$ws.Range("I22").Value = 599.6667
$ws.Range("K22").Value = 599.6667
$ws.Range("M22").Value = -249.6667

$ws.Range("H58").Value = 3995
$ws.Range("J58").Value = 3995
$ws.Range("L58").Value = 3995
$ws.Range("N58").Value = -4401

$ws.Range("H99").Value = 3572.4443
$ws.Range("I99").Value = 3531.5
$ws.Range("K99").Value = 3531.5
$ws.Range("M99").Value = -2033.5

$ws.Range("H105").Value = 995
$ws.Range("I105").Value = 894
$ws.Range("K105").Value = 894
$ws.Range("M105").Value = 853

$ws.Range("H113").Value = 3665
$ws.Range("I113").Value = 1584
$ws.Range("K113").Value = 1584
$ws.Range("M113").Value = 586

$ws.Range("H126").Value = 3572.4443
$ws.Range("I126").Value = 3531.5
$ws.Range("K126").Value = 10594.5
$ws.Range("M126").Value = -8124.5

$ws.Range("H136").Value = 3995
$ws.Range("J136").Value = 3995
$ws.Range("L136").Value = 11985
$ws.Range("N136").Value = -17085

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 11.166667
$ws.Range("I6").Value = 20.333334
$ws.Range("K6").Value = 61.000002
$ws.Range("M6").Value = 51.999998

$ws.Range("H21").Value = 181.25
$ws.Range("I21").Value = 50
$ws.Range("J21").Value = 225
$ws.Range("K21").Value = 150
$ws.Range("L21").Value = 675
$ws.Range("M21").Value = 23
$ws.Range("N21").Value = -1021

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()

$ws.Range("H52").Value = 3950
$ws.Range("J52").Value = 3950
$ws.Range("L52").Value = 11850
$ws.Range("N52").Value = -12382

$ws.Range("H117").Value = 3014.1428
$ws.Range("I117").Value = 804.25
$ws.Range("J117").Value = 5960.6665
$ws.Range("K117").Value = 2412.75
$ws.Range("L117").Value = 17881.9995
$ws.Range("M117").Value = 1029.25
$ws.Range("N117").Value = -24765.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 1100018
$ws.Range("J40").Value = 1100018
$ws.Range("L40").Value = 1100018
$ws.Range("N40").Value = -1100320

$ws.Range("H63").Value = 57229.332
$ws.Range("J63").Value = 57229.332
$ws.Range("L63").Value = 57229.332
$ws.Range("N63").Value = -58601.332

$ws.Range("H66").Value = 57229.332
$ws.Range("J66").Value = 57229.332
$ws.Range("L66").Value = 171687.996
$ws.Range("N66").Value = -178551.996

$ws.Range("H102").Value = 956
$ws.Range("I102").Value = 956
$ws.Range("K102").Value = 956
$ws.Range("M102").Value = 666

$ws.Range("H122").Value = 1434.5555
$ws.Range("I122").Value = 1434.5555
$ws.Range("K122").Value = 4303.666499999999
$ws.Range("M122").Value = -1853.666499999999

$ws.Range("H132").Value = 2194.75
$ws.Range("I132").Value = 1926.6666
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 5779.9998
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -3249.9998
$ws.Range("N132").Value = -14057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 353.8
$ws.Range("J10").Value = 392.5
$ws.Range("L10").Value = 392.5
$ws.Range("N10").Value = -672.5

$ws.Range("H22").Value = 1053.5
$ws.Range("J22").Value = 855.5
$ws.Range("L22").Value = 855.5
$ws.Range("N22").Value = -1445.5

$ws.Range("H27").Value = 1053.5
$ws.Range("J27").Value = 855.5
$ws.Range("L27").Value = 855.5
$ws.Range("N27").Value = -1069.5

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H122").Value = 55569736
$ws.Range("I122").Value = 66682184
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 200046552
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -200044102
$ws.Range("N122").Value = -27400

$ws.Range("H132").Value = 5986.5
$ws.Range("I132").Value = 5986.5
$ws.Range("K132").Value = 17959.5
$ws.Range("M132").Value = -15429.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 2357
$ws.Range("J9").Value = 2357
$ws.Range("L9").Value = 2357
$ws.Range("N9").Value = -2637

$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("M19").ClearContents()

$ws.Range("H51").Value = 39517.75
$ws.Range("I51").Value = 39517.75
$ws.Range("K51").Value = 39517.75
$ws.Range("M51").Value = -39007.75

$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -53494

$ws.Range("H122").Value = 1082.4445
$ws.Range("I122").Value = 1105.375
$ws.Range("K122").Value = 3316.125
$ws.Range("M122").Value = -866.125

$ws.Range("H132").Value = 900
$ws.Range("I132").Value = 900
$ws.Range("K132").Value = 2700
$ws.Range("M132").Value = -170
